$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Acceptance Criteria": update a Gherkin step's text while keeping the
# "Then " keyword bold (matches the rich-text run split in the shared string).
# ---------------------------------------------------------------------------
$wsAC = $wb.Worksheets.Item("Acceptance Criteria")
$cAC = $wsAC.Range("C31")
$cAC.Value = "Then the user should stay on the login page"
$r1 = $cAC.Characters(1, 5)
$r1.Font.Bold = $true
$r1.Font.Name = "Arial"
$r2 = $cAC.Characters(6, 39)
$r2.Font.Bold = $false
$r2.Font.Name = "Arial"

# ---------------------------------------------------------------------------
# Sheet "Test Data": assorted test-data / expected-result text tweaks.
# ---------------------------------------------------------------------------
$wsTD = $wb.Worksheets.Item("Test Data")

# Typo fix in test data value.
$wsTD.Range("Q7").Value = "NotPlanet"

# Tidy up alert copy (remove "alert" suffix / fix casing) for moon & planet name.
# Every cell that shared the old string needs its value re-set explicitly.
$wsTD.Range("N11").Value = "Invalid moon name"
$wsTD.Range("L18").Value = "Invalid moon name"
$wsTD.Range("L19").Value = "Invalid moon name"
$wsTD.Range("L20").Value = "Invalid moon name"
$wsTD.Range("L21").Value = "Invalid moon name"
$wsTD.Range("Q11").Value = "Invalid planet name"

# F14:F17 switch font to Arial (matches style 51 -> 57) and F14's text gains
# detail about the redirect; F15-F17 text drops the "alert" suffix.
$wsTD.Range("F14").Font.Name = "Arial"
$wsTD.Range("F14").Value = "login success, redirect to home"

$wsTD.Range("F15").Font.Name = "Arial"
$wsTD.Range("F15").Value = "Invalid Credentials"
$wsTD.Range("F16").Font.Name = "Arial"
$wsTD.Range("F16").Value = "Invalid Credentials"
$wsTD.Range("F17").Font.Name = "Arial"
$wsTD.Range("F17").Value = "Invalid Credentials"

# I18:I21 previously referenced a duplicate "Invalid planet name alert" string;
# point them at the same cleaned-up text used in Q11 so they dedupe together.
$wsTD.Range("I18").Value = "Invalid planet name"
$wsTD.Range("I19").Value = "Invalid planet name"
$wsTD.Range("I20").Value = "Invalid planet name"
$wsTD.Range("I21").Value = "Invalid planet name"

# I22 / L22 alert copy cleanup.
$wsTD.Range("I22").Value = "Invalid file type"
$wsTD.Range("L22").Value = "Invalid planet ID"

# C24:C36 switch font to Arial (matches style 62 -> new style) and drop the
# "alert" suffix / fix wording on the alert copy.
$wsTD.Range("C24").Font.Name = "Arial"
$wsTD.Range("C24").Value = "Account created successfully"

$wsTD.Range("C25").Font.Name = "Arial"
$wsTD.Range("C25").Value = "Invalid username"
$wsTD.Range("C26").Font.Name = "Arial"
$wsTD.Range("C26").Value = "Invalid username"
$wsTD.Range("C27").Font.Name = "Arial"
$wsTD.Range("C27").Value = "Invalid username"
$wsTD.Range("C28").Font.Name = "Arial"
$wsTD.Range("C28").Value = "Invalid username"
$wsTD.Range("C29").Font.Name = "Arial"
$wsTD.Range("C29").Value = "Invalid username"

$wsTD.Range("C30").Font.Name = "Arial"
$wsTD.Range("C30").Value = "Invalid password"
$wsTD.Range("C31").Font.Name = "Arial"
$wsTD.Range("C31").Value = "Invalid password"
$wsTD.Range("C32").Font.Name = "Arial"
$wsTD.Range("C32").Value = "Invalid password"
$wsTD.Range("C33").Font.Name = "Arial"
$wsTD.Range("C33").Value = "Invalid password"
$wsTD.Range("C34").Font.Name = "Arial"
$wsTD.Range("C34").Value = "Invalid password"
$wsTD.Range("C35").Font.Name = "Arial"
$wsTD.Range("C35").Value = "Invalid password"
$wsTD.Range("C36").Font.Name = "Arial"
$wsTD.Range("C36").Value = "Invalid password"
